$d = $word.ActiveDocument

# The document currently ends with the paragraph containing
# "My name is Prashanth Desai". Add a new paragraph right after it
# with the same paragraph formatting (inherited via InsertParagraphAfter)
# containing the new line of text.
$lastParagraph = $d.Paragraphs($d.Paragraphs.Count)
$lastParagraph.Range.InsertParagraphAfter()

$newParagraph = $d.Paragraphs($d.Paragraphs.Count)
$newParagraph.Range.Text = "i am studying branch command in git"
